# ============================================================================
# Generate Report for Handoff
#
# A new source file, 64ea43da-abd9-4b37-a9e3-77ff1919a36a.md, has reached the
# "Ready for handoff" stage. Insert it as a new row (row 8) ahead of the
# existing 9b2ca098-...md row on every worksheet (Overview, zh-cn, de-de),
# which pushes the 9b2ca098 row and the trailing .localization-config row
# down by one. Hyperlinks are rebuilt from scratch afterwards because the
# Excel row-insert operation does not renumber existing Hyperlink objects.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at position 8. Excel shifts all cell values/styles in
# rows 8+ down by one automatically (the new row also inherits the style
# from the row above it, which matches the formatting we need).
$ws.Rows.Item(8).Insert()

# Populate the newly-opened row 8 with the new files data.
$ws.Range("A8").Value = "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"

# Excel does not keep the Hyperlinks collection in sync with inserted rows,
# so rebuild every hyperlink on the sheet from scratch in the correct order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0c4764f3db3ef2b43709ef899f61f4f61704e2b/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/5b4fac6e-1566-45e6-b942-a94f5b2111ce.md", "", "", "5b4fac6e-1566-45e6-b942-a94f5b2111ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/e2e/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4a58802e30b1a5d6532141ca8813dee6ba4e3ea/e2e/c85cfd81-9223-47a5-b559-a4ac99733a93.md", "", "", "c85cfd81-9223-47a5-b559-a4ac99733a93.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/e079f1c9-825d-4045-a6a3-544c378c7434.md", "", "", "e079f1c9-825d-4045-a6a3-544c378c7434.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/1439b2747a8d4675c6491e1db292a233079da5d5/e2e/64ea43da-abd9-4b37-a9e3-77ff1919a36a.md", "", "", "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/71f0e5d7b1b1511f60c5dac8580283087f045948/e2e/9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md", "", "", "9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/.localization-config", "", "", ".localization-config") | Out-Null

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Insert a blank row at position 8. Excel shifts all cell values/styles in
# rows 8+ down by one automatically (the new row also inherits the style
# from the row above it, which matches the formatting we need).
$ws.Rows.Item(8).Insert()

# Populate the newly-opened row 8 with the new files data.
$ws.Range("A8").Value = "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.zh-cn.xlf"
$ws.Range("D8").Value = "2016-03-10 18:41:39"
$ws.Range("G8").Value = "0001-01-01 00:00:00"
$ws.Range("H8").Value = "Include"

# Excel does not keep the Hyperlinks collection in sync with inserted rows,
# so rebuild every hyperlink on the sheet from scratch in the correct order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0c4764f3db3ef2b43709ef899f61f4f61704e2b/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e65636c73c39842a997b399f179f315bccf88e89/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.zh-cn.xlf", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3a32046c6e73beb812164415d6c39dad530d6f8b/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/048bd063ff5b9f12cd087220a794c155a40a0e42/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.zh-cn.xlf", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/5b4fac6e-1566-45e6-b942-a94f5b2111ce.md", "", "", "5b4fac6e-1566-45e6-b942-a94f5b2111ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec525909006d271fc3d3f3bc0cd39db1a8a582b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5b4fac6e-1566-45e6-b942-a94f5b2111ce.60e93ca4dd9ac7f1f1137df167e0876f2c16d9a9.zh-cn.xlf", "", "", "5b4fac6e-1566-45e6-b942-a94f5b2111ce.60e93ca4dd9ac7f1f1137df167e0876f2c16d9a9.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/e2e/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4bb90de115b9f3203bc3a14ffa54950b6c31999/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.zh-cn.xlf", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c2bcae88eee29899cb8abf734cf7dd0b91589aa4/e2e/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/25a8e6dc4a7195e8999beef5d9d59a2cec5532df/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.zh-cn.xlf", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4a58802e30b1a5d6532141ca8813dee6ba4e3ea/e2e/c85cfd81-9223-47a5-b559-a4ac99733a93.md", "", "", "c85cfd81-9223-47a5-b559-a4ac99733a93.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94fc68092a8cf36adcb5626b754db89e8c0a53dc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85cfd81-9223-47a5-b559-a4ac99733a93.2f0553c99450b2553df3200ab82f250f772d1564.zh-cn.xlf", "", "", "c85cfd81-9223-47a5-b559-a4ac99733a93.2f0553c99450b2553df3200ab82f250f772d1564.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/e079f1c9-825d-4045-a6a3-544c378c7434.md", "", "", "e079f1c9-825d-4045-a6a3-544c378c7434.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec525909006d271fc3d3f3bc0cd39db1a8a582b3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e079f1c9-825d-4045-a6a3-544c378c7434.dec076c3cd421d849d05efe05d476ba8fdb7b355.zh-cn.xlf", "", "", "e079f1c9-825d-4045-a6a3-544c378c7434.dec076c3cd421d849d05efe05d476ba8fdb7b355.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cab619e0fff4d4bcaf590b1bfab8fd612acb9ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.zh-cn.xlf", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/1439b2747a8d4675c6491e1db292a233079da5d5/e2e/64ea43da-abd9-4b37-a9e3-77ff1919a36a.md", "", "", "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9cbe2b814d5fe750e8d22be6a5492fba930954b1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.zh-cn.xlf", "", "", "64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/71f0e5d7b1b1511f60c5dac8580283087f045948/e2e/9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md", "", "", "9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e880af9c500f6b1e4e98b74610c5c7bb57df78b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9b2ca098-1e00-441b-aa9b-d3965ac92bb8.8fabbd82688d7ff9fb19711da05e23b803c9a2ce.zh-cn.xlf", "", "", "9b2ca098-1e00-441b-aa9b-d3965ac92bb8.8fabbd82688d7ff9fb19711da05e23b803c9a2ce.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/.localization-config", "", "", ".localization-config") | Out-Null

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

# Insert a blank row at position 8. Excel shifts all cell values/styles in
# rows 8+ down by one automatically (the new row also inherits the style
# from the row above it, which matches the formatting we need).
$ws.Rows.Item(8).Insert()

# Populate the newly-opened row 8 with the new files data.
$ws.Range("A8").Value = "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.de-de.xlf"
$ws.Range("D8").Value = "2016-03-10 18:41:44"
$ws.Range("G8").Value = "0001-01-01 00:00:00"
$ws.Range("H8").Value = "Include"

# Excel does not keep the Hyperlinks collection in sync with inserted rows,
# so rebuild every hyperlink on the sheet from scratch in the correct order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c0c4764f3db3ef2b43709ef899f61f4f61704e2b/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b79ba79e982f27610a82da80f7982ec2e38f6bb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.de-de.xlf", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e53afca1cac254205ccf0e0525c39fb118d7f640/e2e/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/182c656a21acb49e01e496981cdb830866f4055e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.de-de.xlf", "", "", "bf64de53-ac12-4f98-b4f7-d6e6011dc4b1.8438a1e037094bca3c3cb27aaec2309b95e6b8d9.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/5b4fac6e-1566-45e6-b942-a94f5b2111ce.md", "", "", "5b4fac6e-1566-45e6-b942-a94f5b2111ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dca339778b041a6ae151beb0427194c311769a66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5b4fac6e-1566-45e6-b942-a94f5b2111ce.60e93ca4dd9ac7f1f1137df167e0876f2c16d9a9.de-de.xlf", "", "", "5b4fac6e-1566-45e6-b942-a94f5b2111ce.60e93ca4dd9ac7f1f1137df167e0876f2c16d9a9.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/e2e/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/446dca3ae1519f8979ca4c9ac3a29432bed223cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.de-de.xlf", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/18cb633c24ddc9e8ea8eb3c7c9c9de3078d9301b/e2e/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b30dc190668a971a4c779d879859ccdefc77ce35/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.de-de.xlf", "", "", "88dcc02c-743f-49e0-9bdb-ee929a4ebee4.b3474cae98f0230cdf3700f8ea440bc82a708717.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f4a58802e30b1a5d6532141ca8813dee6ba4e3ea/e2e/c85cfd81-9223-47a5-b559-a4ac99733a93.md", "", "", "c85cfd81-9223-47a5-b559-a4ac99733a93.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f295fe709ae672a9d032ee6512209c6066d0c827/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85cfd81-9223-47a5-b559-a4ac99733a93.2f0553c99450b2553df3200ab82f250f772d1564.de-de.xlf", "", "", "c85cfd81-9223-47a5-b559-a4ac99733a93.2f0553c99450b2553df3200ab82f250f772d1564.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a87afa9fd205763fd7be909fe20214c90cb037f4/e2e/e079f1c9-825d-4045-a6a3-544c378c7434.md", "", "", "e079f1c9-825d-4045-a6a3-544c378c7434.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dca339778b041a6ae151beb0427194c311769a66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e079f1c9-825d-4045-a6a3-544c378c7434.dec076c3cd421d849d05efe05d476ba8fdb7b355.de-de.xlf", "", "", "e079f1c9-825d-4045-a6a3-544c378c7434.dec076c3cd421d849d05efe05d476ba8fdb7b355.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c0930bafda350e20443d571f280161b768437c6d/e2e/00cbf9ee-0a9a-4a30-84ad-1ec060825999.md", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8bbebc01928feaf0b4b60fad591eee00c9518ad4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.de-de.xlf", "", "", "00cbf9ee-0a9a-4a30-84ad-1ec060825999.a0e2c0d7e0fed157e66419b3476280ec3bdc3c13.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/1439b2747a8d4675c6491e1db292a233079da5d5/e2e/64ea43da-abd9-4b37-a9e3-77ff1919a36a.md", "", "", "64ea43da-abd9-4b37-a9e3-77ff1919a36a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b6f2e9456d3ad8f2590dc6435e80a4761a125f1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.de-de.xlf", "", "", "64ea43da-abd9-4b37-a9e3-77ff1919a36a.9ac4d004126dbf7295b70e29a2467630219b250c.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/71f0e5d7b1b1511f60c5dac8580283087f045948/e2e/9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md", "", "", "9b2ca098-1e00-441b-aa9b-d3965ac92bb8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89cf1caa248a4bf9afa80e56b396344c7a06a3f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9b2ca098-1e00-441b-aa9b-d3965ac92bb8.8fabbd82688d7ff9fb19711da05e23b803c9a2ce.de-de.xlf", "", "", "9b2ca098-1e00-441b-aa9b-d3965ac92bb8.8fabbd82688d7ff9fb19711da05e23b803c9a2ce.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/6113da6f937d301f22b6d8da7bbcb6f389d9f5e3/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report for Handoff generated: inserted 64ea43da-abd9-4b37-a9e3-77ff1919a36a.md on all sheets."
